$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.287.66"
$ws.Range("E2").Value = "  -4.22%  "

$ws.Range("D3").Value = "3.115.51"
$ws.Range("E3").Value = "  -4.31%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'605.64"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("D6").Value = "'144.31"
$ws.Range("E6").Value = "  -8.21%  "

$ws.Range("D8").Value = "3.112.75"
$ws.Range("E8").Value = "  -4.36%  "

$ws.Range("E9").Value = "  -4.42%  "

$ws.Range("E10").Value = "  -6.97%  "

$ws.Range("D11").Value = "'5.22"
$ws.Range("E11").Value = "  -8.47%  "

$ws.Range("E12").Value = "  -5.60%  "

$ws.Range("D13").Value = "'0.0000248"
$ws.Range("E13").Value = "  -6.81%  "

$ws.Range("D14").Value = "'35.02"
$ws.Range("E14").Value = "  -9.12%  "

$ws.Range("D15").Value = "3.614.85"
$ws.Range("E15").Value = "  -4.35%  "

$ws.Range("E16").Value = "  +1.49%  "

$ws.Range("D17").Value = "63.391.18"
$ws.Range("E17").Value = "  -4.04%  "

$ws.Range("D18").Value = "3.104.83"
$ws.Range("E18").Value = "  -4.48%  "

$ws.Range("D19").Value = "'6.76"
$ws.Range("E19").Value = "  -7.23%  "

$ws.Range("D20").Value = "'471.86"
$ws.Range("E20").Value = "  -5.23%  "

$ws.Range("D21").Value = "'14.48"
$ws.Range("E21").Value = "  -5.27%  "

$ws.Range("D22").Value = "'0.701"
$ws.Range("E22").Value = "  -6.06%  "

$ws.Range("D23").Value = "'7.67"
$ws.Range("E23").Value = "  -4.44%  "

$ws.Range("D24").Value = "'13.44"
$ws.Range("E24").Value = "  -7.77%  "

$ws.Range("D25").Value = "'82.73"
$ws.Range("E25").Value = "  -4.84%  "

$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("E27").Value = "  -8.44%  "

$ws.Range("E28").Value = "  -8.06%  "

$ws.Range("D29").Value = "'0.119"
$ws.Range("E29").Value = "  -9.83%  "

$ws.Range("D30").Value = "'6.83"
$ws.Range("E30").Value = "  -3.28%  "

$ws.Range("D31").Value = "'2.07"
$ws.Range("E31").Value = "  -12.38%  "

$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").Value = "'2.65"
$ws.Range("E33").Value = "  -6.78%  "

$ws.Range("D34").Value = "'26.11"
$ws.Range("E34").Value = "  -6.18%  "

$ws.Range("D35").Value = "'1.10"
$ws.Range("E35").Value = "  -2.81%  "

$ws.Range("D36").Value = "'5.89"
$ws.Range("E36").Value = "  -7.66%  "

$ws.Range("D37").Value = "'52.61"
$ws.Range("E37").Value = "  -5.48%  "

$ws.Range("D38").Value = "0.0₃0749"
$ws.Range("E38").Value = "  -2.42%  "

$ws.Range("D39").Value = "'452.86"
$ws.Range("E39").Value = "  -8.18%  "

$ws.Range("E40").Value = "  -13.91%  "

$ws.Range("D41").Value = "'0.0390"
$ws.Range("E41").Value = "  -7.19%  "

$ws.Range("E42").Value = "  -9.67%  "

$ws.Range("E43").Value = "  -5.47%  "

$ws.Range("D44").Value = "2.832.02"
$ws.Range("E44").Value = "  -5.41%  "

$ws.Range("E45").Value = "  -10.79%  "

$ws.Range("D46").Value = "'0.262"
$ws.Range("E46").Value = "  -9.82%  "

$ws.Range("E47").Value = "  -2.29%  "

$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").Value = "'25.95"
$ws.Range("E49").Value = "  -9.25%  "

$ws.Range("E50").Value = "  -5.76%  "

$ws.Range("D51").Value = "'118.86"
$ws.Range("E51").Value = "  -1.96%  "
